# Adds the new "Asset" module/page/permission/ability rows (sibling to the
# existing "Infrastructure Unit Asset" rows, but without the "Unit" level)
# across the pages / permissions / abilities-pages / abilities-permissions
# sheets, and updates the active-sheet/selection bookkeeping to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "pages": new row 5 - the "Asset" page entry
# ---------------------------------------------------------------------
$pages = $wb.Worksheets.Item("pages")

$pages.Range("A5").Value = "infrastructure"
$pages.Range("B5").Value = "Infrastructure Asset"
$pages.Range("C5").Value = "Asset"
$pages.Range("D5").Value = "infrastructure-asset"
$pages.Range("E5").Value = "tune"
$pages.Range("F5").Value = "infrastructure"
$pages.Range("G5").Value = "asset"
$pages.Range("H5").Formula = "=TRUE()"
$pages.Range("I5").Formula = "=FALSE()"
$pages.Range("K5").Formula = "=TRUE()"

$pages.Activate()
$pages.Range("C4").Select()

# ---------------------------------------------------------------------
# Sheet "permissions": new row 5 - permissions for the "Asset" module
# ---------------------------------------------------------------------
$permissions = $wb.Worksheets.Item("permissions")

$permissions.Range("A5").Value = "infrastructure"
$permissions.Range("B5").Value = "infrastructure-asset"
$permissions.Range("C5").Value = "view, create, show, update, delete, restore, destroy"

$permissions.Activate()
$permissions.Range("B5").Select()

# ---------------------------------------------------------------------
# Sheet "abilities-pages": new row 5 - ability/page link for "Asset"
# ---------------------------------------------------------------------
$abilitiesPages = $wb.Worksheets.Item("abilities-pages")

$abilitiesPages.Range("A5").Value = "infrastructure"
$abilitiesPages.Range("B5").Value = "infrastructure-asset"
$abilitiesPages.Range("C5").Value = "superadmin, administrator"

$abilitiesPages.Activate()
$abilitiesPages.Range("B5").Select()

# ---------------------------------------------------------------------
# Sheet "abilities-permissions": new row 5 - ability/permission link for "Asset"
# ---------------------------------------------------------------------
$abilitiesPermissions = $wb.Worksheets.Item("abilities-permissions")

$abilitiesPermissions.Range("A5").Value = "infrastructure"
$abilitiesPermissions.Range("B5").Value = "infrastructure-asset"
$abilitiesPermissions.Range("C5").Value = "superadmin, administrator"
$abilitiesPermissions.Range("D5").Value = "*"

# This is the sheet that ends up active/selected when the workbook is saved
# (workbook.xml bookViews activeTab moves from 1 -> 5).
$abilitiesPermissions.Activate()
$abilitiesPermissions.Range("B5").Select()
